$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a "plain number" string must be forced back to
# Text format first, otherwise Excel silently re-types them as numbers
# (and can even drop a significant trailing zero, e.g. "64.60" -> 64.6).
# The source workbook stores every Price/Volume cell as text, so we keep
# that invariant for every numeric-looking replacement value.

$ws.Range("D2").Value = '29.863.25'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.623.00'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.37'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.74'
$ws.Range("E8").Value = '  +10.59%  '
$ws.Range("E9").Value = '  +2.77%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '1.855.62'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").Value = '1.624.47'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("E14").Value = '  +5.81%  '
$ws.Range("E15").Value = '  +4.75%  '
$ws.Range("D16").Value = '29.922.43'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.81'
$ws.Range("E17").Value = '  +16.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.60'
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.68'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.11'
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.60'
$ws.Range("E23").Value = '  +4.04%  '
$ws.Range("E24").Value = '  +2.38%  '
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.67'
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.111'
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("E28").Value = '  +2.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0488'
$ws.Range("E30").Value = '  +3.09%  '
$ws.Range("E31").Value = '  +5.15%  '
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("E33").Value = '  +3.40%  '
$ws.Range("D34").Value = '1.427.54'
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("E35").Value = '  +6.91%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("E39").Value = '  +2.98%  '
$ws.Range("E40").Value = '  +3.32%  '
$ws.Range("E41").Value = '  +3.35%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("E43").Value = '  +4.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '53.96'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.18'
$ws.Range("E46").Value = '  +17.74%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("E48").Value = '  +2.49%  '
$ws.Range("D49").Value = '1.764.10'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.39'
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("D51").Value = '0.0₆0106'
$ws.Range("E51").Value = '  +1.00%  '
